# Update the "dSF" column (F) values for several rows.
# These edits correspond to a repull/recalculation of the data
# (see commit message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    9  = 3
    11 = -3
    13 = 2
    18 = 8
    19 = -1
    21 = -2
    27 = -1
    30 = -2
    31 = -4
    35 = 1
    37 = 8
    44 = -1
    47 = -3
    48 = -3
    50 = -8
    52 = -2
    54 = -3
    57 = 1
    63 = 2
    66 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
